# Update automatico via Actualizar 07-29-2020 00-42-56
# Adds 4 new rows (30-33) of MITRADEL news items to the "trabajo" table,
# plus one new trailing blank table row (34), expands the table /
# autofilter / data validation ranges accordingly, wires up the new
# hyperlinks, and fixes up row heights + the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---------------------------------------------------------------------
# 1. Grow the table so the new rows inherit the table's formatting/
#    autofilter/etc. Row 34 is the fresh blank row at the bottom that a
#    table gets when its range is extended without being filled in.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K34"))

# ---------------------------------------------------------------------
# 2. Seed formatting for the new rows by copying from existing rows that
#    already carry the exact styles the new rows need.
#    Row 27 -> rows 30, 31, 34 (plain "Ministerio de Trabajo" row style)
# ---------------------------------------------------------------------
$ws.Range("A27:K27").Copy() | Out-Null
$ws.Range("A30:K30").PasteSpecial(-4122) | Out-Null
$ws.Range("A31:K31").PasteSpecial(-4122) | Out-Null
$ws.Range("A34:K34").PasteSpecial(-4122) | Out-Null

# Rows 32, 33 use the "header-style-2" pattern seen in rows 2-10 / 21,
# built column by column since no single existing row matches exactly.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A32").PasteSpecial(-4122) | Out-Null
$ws.Range("A33").PasteSpecial(-4122) | Out-Null

$ws.Range("B21").Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4122) | Out-Null
$ws.Range("B33").PasteSpecial(-4122) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("C32").PasteSpecial(-4122) | Out-Null
$ws.Range("C33").PasteSpecial(-4122) | Out-Null

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4122) | Out-Null
$ws.Range("D33").PasteSpecial(-4122) | Out-Null

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E32").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").PasteSpecial(-4122) | Out-Null

$ws.Range("F2").Copy() | Out-Null
$ws.Range("F32").PasteSpecial(-4122) | Out-Null
$ws.Range("F33").PasteSpecial(-4122) | Out-Null

$ws.Range("G21").Copy() | Out-Null
$ws.Range("G32").PasteSpecial(-4122) | Out-Null
$ws.Range("G33").PasteSpecial(-4122) | Out-Null

$ws.Range("H2").Copy() | Out-Null
$ws.Range("H32").PasteSpecial(-4122) | Out-Null
$ws.Range("H33").PasteSpecial(-4122) | Out-Null

$ws.Range("I2").Copy() | Out-Null
$ws.Range("I32").PasteSpecial(-4122) | Out-Null
$ws.Range("I33").PasteSpecial(-4122) | Out-Null

$ws.Range("J2").Copy() | Out-Null
$ws.Range("J32").PasteSpecial(-4122) | Out-Null
$ws.Range("J33").PasteSpecial(-4122) | Out-Null

$ws.Range("K2").Copy() | Out-Null
$ws.Range("K32").PasteSpecial(-4122) | Out-Null
$ws.Range("K33").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Row heights
# ---------------------------------------------------------------------
$ws.Rows(30).RowHeight = 120
$ws.Rows(31).RowHeight = 140.25
$ws.Rows(32).RowHeight = 120
$ws.Rows(33).RowHeight = 120
$ws.Rows(34).RowHeight = 18.75

# ---------------------------------------------------------------------
# 4. Cell content
# ---------------------------------------------------------------------

# -- Row 30 : Proyecto de Ley de proteccion del empleo (anuncio) --
$ws.Range("B30").Value = 29
$e30 = @"
https://www.mitradel.gob.pa/la-proxima-semana-se-presentara-a-la-asamblea-nacional-el-proyecto-de-ley-de-proteccion-del-empleo-en-empresas-afectadas-por-el-covid-19/
"@
$ws.Range("E30").Value = $e30
$f30 = @"
El próximo lunes 13 de julio de 2020, la Ministra Doris Zapata Acevedo presentará a la Asamblea Nacional el Proyecto de Ley que establece medidas temporales de protección del empleo, aplicable únicamente en las empresas que cerraron, total o parcialmente sus operaciones, desde el inicio del Estado de Emergencia Nacional y que fueron consideradas en la Mesa Tripartita Económico Laboral.
"@
$ws.Range("F30").Value = $f30
$ws.Range("H30").Value = "28-07-2020"
$ws.Range("I30").Value = 44022

# -- Row 31 : Modificacion de la jornada de trabajo --
$ws.Range("B31").Value = 30
$e31 = @"
https://www.mitradel.gob.pa/modificacion-de-la-jornada-de-trabajo-tienen-regulaciones/
"@
$ws.Range("E31").Value = $e31
$f31 = @"
A través del Decreto Ejecutivo No.101 del 13 de julio de 2020, el Ministerio de Trabajo y Desarrollo Laboral formalizó la norma para la convención de la modificación o reducción temporal de la jornada de trabajo, en cumplimiento de uno de los acuerdos alcanzados en la Mesa Tripartita de Diálogo por la Economía y el Desarrollo Laboral.
La nueva norma establece que en el acuerdo de modificación de la jornada laboral se deben incluir métodos para lograr la recuperación gradual de las jornadas laborales a los niveles existentes antes de la crisis y que no deberán afectar la rata por hora pactada en el contrato de trabajo vigente.
"@
$ws.Range("F31").Value = $f31
$ws.Range("H31").Value = "28-07-2020"
$ws.Range("I31").Value = 44025

# -- Row 32 : Proponen ley para ampliar vale de alimentacion --
$ws.Range("A32").Value = "Ministerio de Trabajo y Desarrollo Social"
$ws.Range("B32").Value = 31
$ws.Range("C32").Value = "Trabajo"
$ws.Range("D32").Value = "El Ministerio de Trabajo y Desarrollo Laboral de Panamá (MITRADEL) es un Ministerio de la República de Panamá que forma parte del Órgano Ejecutivo, encargado de regular las relaciones obrero-patronales, fomentar la generación de empleo y la formación de la mano de obra nacional, así como de velar por la aplicación de las normas laborales."
$e32 = @"
https://www.mitradel.gob.pa/proponen-ley-para-ampliar-el-alcance-del-vale-de-alimentacion/
"@
$ws.Range("E32").Value = $e32
$f32 = @"
En cumplimiento de uno de los acuerdos alcanzados en la Mesa Tripartita de Diálogo por la Economía y el Desarrollo Laboral, la Ministra de Trabajo y Desarrollo Laboral, Doris Zapata Acevedo, presentó este martes 21 de julio de 2020, ante el Consejo de Gabinete, un Proyecto de Ley que contiene la propuesta de modificación al Programa de Alimentación para los Trabajadores, como una medida orientada a ampliar los beneficios otorgados a través de este programa.
"@
$ws.Range("F32").Value = $f32
$ws.Range("G32").Value = "https://www.mitradel.gob.pa"
$ws.Range("H32").Value = "28-07-2020"
$ws.Range("I32").Value = "21-07-2020"
$ws.Range("J32").Value = "Panamá"
$ws.Range("K32").Value = "Ministerial"

# -- Row 33 : Inicia primer debate del proyecto de ley --
$ws.Range("A33").Value = "Ministerio de Trabajo y Desarrollo Social"
$ws.Range("B33").Value = 32
$ws.Range("C33").Value = "Trabajo"
$ws.Range("D33").Value = "El Ministerio de Trabajo y Desarrollo Laboral de Panamá (MITRADEL) es un Ministerio de la República de Panamá que forma parte del Órgano Ejecutivo, encargado de regular las relaciones obrero-patronales, fomentar la generación de empleo y la formación de la mano de obra nacional, así como de velar por la aplicación de las normas laborales."
$e33 = @"
https://www.mitradel.gob.pa/inicia-el-primer-debate-del-proyecto-de-ley-de-proteccion-del-empleo/
"@
$ws.Range("E33").Value = $e33
$f33 = @"
La Comisión de Salud, Trabajo y Desarrollo Social de la Asamblea Nacional inició, este jueves 23 de julio de 2020, el Primer Debate del Proyecto de Ley 354 que establece medidas temporales de protección del empleo, aplicable únicamente en las empresas que cerraron, total o parcialmente, sus operaciones, desde el inicio del Estado de Emergencia Nacional.
"@
$ws.Range("F33").Value = $f33
$ws.Range("G33").Value = "https://www.mitradel.gob.pa"
$ws.Range("H33").Value = "28-07-2020"
$ws.Range("I33").Value = "23-07-2020"
$ws.Range("J33").Value = "Panamá"
$ws.Range("K33").Value = "Ministerial"

# Row 34 stays blank (just formatted) - matches the fresh trailing table
# row Excel leaves when the table is extended past the last filled row.

# ---------------------------------------------------------------------
# 5. Hyperlinks for the new cells
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E30"), $e30.Trim()) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E31"), $e31.Trim()) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E32"), $e32.Trim()) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E33"), $e33.Trim()) | Out-Null
$ws.Hyperlinks.Add($ws.Range("G32"), "https://www.mitradel.gob.pa/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G33"), "https://www.mitradel.gob.pa/") | Out-Null

# ---------------------------------------------------------------------
# 6. Data validation range grows along with the table (C2:C34)
# ---------------------------------------------------------------------
$ws.Range("C2:C34").Validation.Delete()
$ws.Range("C2:C34").Validation.Add(0)
$ws.Range("C2:C34").Validation.IgnoreBlank = $true
$ws.Range("C2:C34").Validation.InCellDropdown = $true
$ws.Range("C2:C34").Validation.ErrorTitle = "Entrada no válida"
$ws.Range("C2:C34").Validation.ErrorMessage = "Selecciona una categoría de la lista"
$ws.Range("C2:C34").Validation.InputTitle = "Categoria"
$ws.Range("C2:C34").Validation.InputMessage = "Selecciona una categoría de la lista"

# ---------------------------------------------------------------------
# 7. Selection / view state
# ---------------------------------------------------------------------
$ws.Range("H30:H33").Select()
